# Applies odds updates to rows 3, 10, 11, 12, and 14 of the active sheet,
# matching the target diff for Jogos_da_Semana_FlashScore_2025-05-27.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$ws.Range("H3").Value = 3.05
$ws.Range("I3").Value = 4.15
$ws.Range("L3").Value = 1.42
$ws.Range("M3").Value = 2.47
$ws.Range("N3").Value = 2.22
$ws.Range("O3").Value = 1.52
$ws.Range("P3").Value = 1.47
$ws.Range("Q3").Value = 2.32
$ws.Range("R3").Value = 1.98
$ws.Range("S3").Value = 1.65
$ws.Range("T3").Value = 5.5
$ws.Range("U3").Value = 7.9
$ws.Range("V3").Value = 8.75
$ws.Range("W3").Value = 16.5
$ws.Range("X3").Value = 18.5
$ws.Range("Y3").Value = 37
$ws.Range("Z3").Value = 7.1
$ws.Range("AA3").Value = 6
$ws.Range("AB3").Value = 17
$ws.Range("AC3").Value = 100
$ws.Range("AE3").Value = 9.75
$ws.Range("AG3").Value = 14
$ws.Range("AH3").Value = 70
$ws.Range("AI3").Value = 45
$ws.Range("AJ3").Value = 55

# --- Row 10 ---
$ws.Range("N10").Value = 1.89
$ws.Range("O10").Value = 1.79
$ws.Range("R10").Value = 1.67
$ws.Range("S10").Value = 2.1
$ws.Range("T10").Value = 8.5
$ws.Range("U10").Value = 12
$ws.Range("X10").Value = 17
$ws.Range("Y10").Value = 26
$ws.Range("Z10").Value = 11
$ws.Range("AB10").Value = 13
$ws.Range("AC10").Value = 41
$ws.Range("AD10").Value = 201
$ws.Range("AE10").Value = 10
$ws.Range("AI10").Value = 23
$ws.Range("AJ10").Value = 29

# --- Row 11 ---
$ws.Range("L11").Value = 1.29
$ws.Range("M11").Value = 3.5
$ws.Range("N11").Value = 1.98
$ws.Range("O11").Value = 1.83

# --- Row 12 ---
$ws.Range("G12").Value = 1.83
$ws.Range("I12").Value = 3.8
$ws.Range("U12").Value = 9
$ws.Range("W12").Value = 15
$ws.Range("AA12").Value = 7.5
$ws.Range("AC12").Value = 51

# --- Row 14 ---
$ws.Range("G14").Value = 2
$ws.Range("I14").Value = 3.4
$ws.Range("N14").Value = 1.67
$ws.Range("O14").Value = 2.15
$ws.Range("P14").Value = 1.3
$ws.Range("Q14").Value = 3.4
$ws.Range("Y14").Value = 21
$ws.Range("AH14").Value = 41
